# Added two new Mac-Addresses (10 new device rows, rows 147-156)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$startRow = 147
$startDevice = 3000166

for ($i = 0; $i -lt 10; $i++) {
    $row = $startRow + $i
    $deviceId = $startDevice + $i

    $ws.Cells.Item($row, 1).Value = 10001
    $ws.Cells.Item($row, 2).Value = $deviceId
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
    $ws.Cells.Item($row, 7).Value = "now()"
}

# Update the selected/active cell to match the new end of the data range
$ws.Range("E155").Select()
